$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E, rows 16-32) used to list the overdue periods
# in ascending order (2209 .. 2401). The update removes the old periods and
# adds the new ones, effectively listing them in descending (most recent
# first) order: 2401 .. 2209.
$periods = @("2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212","2211","2210","2209")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}
